$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from existing header cell (H1) to new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for columns I (I0) and J (IF)
$values = @{
    2  = @(10, 11)
    3  = @(6, 6)
    4  = @(7, 7)
    5  = @(8, 8)
    6  = @(7, 8)
    7  = @(1, 1)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(8, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
